$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue 'D2' '43.276.73'
Set-TextValue 'D3' '2.279.43'
Set-TextValue 'E3' '  -1.17%  '
Set-TextValue 'E4' '  +0.01%  '
Set-TextValue 'D5' '112.48'
Set-TextValue 'E5' '  -1.98%  '
Set-TextValue 'D6' '265.22'
Set-TextValue 'E6' '  -1.57%  '
Set-TextValue 'E7' '  -0.87%  '
Set-TextValue 'D9' '0.609'
Set-TextValue 'E9' '  -2.77%  '
Set-TextValue 'D10' '47.51'
Set-TextValue 'E10' '  -2.21%  '
Set-TextValue 'D11' '0.0929'
Set-TextValue 'E11' '  -1.67%  '
Set-TextValue 'E12' '  +2.03%  '
Set-TextValue 'E13' '  +0.93%  '
Set-TextValue 'D14' '15.52'
Set-TextValue 'E14' '  -1.11%  '
Set-TextValue 'D15' '2.622.73'
Set-TextValue 'E15' '  -0.79%  '
Set-TextValue 'D16' '0.859'
Set-TextValue 'E16' '  -0.46%  '
Set-TextValue 'D17' '2.280.22'
Set-TextValue 'E17' '  -1.34%  '
Set-TextValue 'D18' '43.168.50'
Set-TextValue 'E18' '  -1.33%  '
Set-TextValue 'E19' '  -2.52%  '
Set-TextValue 'E20' '  +3.04%  '
Set-TextValue 'D21' '71.44'
Set-TextValue 'E21' '  -1.68%  '
Set-TextValue 'E22' '  -2.33%  '
Set-TextValue 'D23' '232.07'
Set-TextValue 'E23' '  -1.13%  '
Set-TextValue 'D24' '9.62'
Set-TextValue 'E24' '  +0.75%  '
Set-TextValue 'E25' '  +0.13%  '
Set-TextValue 'E26' '  +0.97%  '
Set-TextValue 'D27' '11.31'
Set-TextValue 'E27' '  -2.27%  '
Set-TextValue 'B28' 'LEO'
Set-TextValue 'C28' 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue 'D28' '3.92'
Set-TextValue 'E28' '  -0.93%  '
Set-TextValue 'B29' 'InjectiveProtocol'
Set-TextValue 'C29' 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue 'D29' '40.42'
Set-TextValue 'E29' '  -7.44%  '
Set-TextValue 'B30' 'WEMIXToken'
Set-TextValue 'C30' 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue 'D30' '3.35'
Set-TextValue 'E30' '  -1.98%  '
Set-TextValue 'B31' 'Toncoin'
Set-TextValue 'C31' 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue 'D31' '2.25'
Set-TextValue 'E31' '  -1.11%  '
Set-TextValue 'B32' 'Monero'
Set-TextValue 'C32' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 'D32' '171.89'
Set-TextValue 'E32' '  -3.46%  '
Set-TextValue 'B33' 'EthereumClassic'
Set-TextValue 'C33' 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue 'D33' '21.34'
Set-TextValue 'E33' '  -2.72%  '
Set-TextValue 'B34' 'Hedera'
Set-TextValue 'C34' 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 'D34' '0.0907'
Set-TextValue 'E34' '  -3.04%  '
Set-TextValue 'B35' 'Filecoin'
Set-TextValue 'C35' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D35' '5.80'
Set-TextValue 'E35' '  +4.37%  '
Set-TextValue 'B36' 'Stellar'
Set-TextValue 'C36' 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue 'D36' '0.127'
Set-TextValue 'E36' '  +0.25%  '
Set-TextValue 'B37' 'RenderToken'
Set-TextValue 'C37' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D37' '4.63'
Set-TextValue 'E37' '  -2.49%  '
Set-TextValue 'B38' 'NEARProtocol'
Set-TextValue 'C38' 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue 'D38' '3.92'
Set-TextValue 'E38' '  -1.03%  '
Set-TextValue 'B39' 'VeChain'
Set-TextValue 'C39' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D39' '0.0355'
Set-TextValue 'E39' '  -0.59%  '
Set-TextValue 'B40' 'Kaspa'
Set-TextValue 'C40' 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue 'D40' '0.104'
Set-TextValue 'E40' '  -5.54%  '
Set-TextValue 'B41' 'LidoDAOToken'
Set-TextValue 'C41' 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue 'D41' '2.63'
Set-TextValue 'E41' '  +9.22%  '
Set-TextValue 'B42' 'MultiversX'
Set-TextValue 'C42' 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
Set-TextValue 'D42' '79.13'
Set-TextValue 'E42' '  +4.65%  '
Set-TextValue 'B43' 'Celestia'
Set-TextValue 'C43' 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
Set-TextValue 'D43' '13.81'
Set-TextValue 'E43' '  +3.16%  '
Set-TextValue 'B44' 'Algorand'
Set-TextValue 'C44' 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue 'D44' '0.238'
Set-TextValue 'E44' '  -4.04%  '
Set-TextValue 'B45' 'THORChain'
Set-TextValue 'C45' 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
Set-TextValue 'D45' '6.25'
Set-TextValue 'E45' '  +3.97%  '
Set-TextValue 'B46' 'FirstDigitalUSD'
Set-TextValue 'C46' 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue 'D46' '1.00'
Set-TextValue 'E46' '  -0.24%  '
Set-TextValue 'B47' 'ARBITRUM'
Set-TextValue 'C47' 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue 'D47' '1.39'
Set-TextValue 'E47' '  -2.62%  '
Set-TextValue 'B48' 'FraxShare'
Set-TextValue 'C48' 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue 'D48' '8.69'
Set-TextValue 'E48' '  -0.82%  '
Set-TextValue 'B49' 'Aave'
Set-TextValue 'C49' 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 'D49' '103.58'
Set-TextValue 'E49' '  +1.68%  '
Set-TextValue 'B50' 'TrustWalletToken'
Set-TextValue 'C50' 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue 'D50' '1.25'
Set-TextValue 'E50' '  +1.45%  '
Set-TextValue 'B51' 'Cronos'
Set-TextValue 'C51' 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue 'D51' '0.0994'
Set-TextValue 'E51' '  -1.47%  '
